$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 63

# Row 15
$ws.Range("E15").Value = 181
$ws.Range("F15").Value = 101
$ws.Range("H15").Value = 142

# Row 17
$ws.Range("E17").Value = 142
$ws.Range("F17").Value = 74
$ws.Range("H17").Value = 106

# Row 19
$ws.Range("E19").Value = 71

# Row 23
$ws.Range("E23").Value = 9

# Row 26
$ws.Range("E26").Value = 36
$ws.Range("F26").Value = 20
$ws.Range("H26").Value = 30

# Row 28
$ws.Range("E28").Value = 23
$ws.Range("F28").Value = 20
$ws.Range("H28").Value = 22

# Row 32
$ws.Range("E32").Value = 24

# Row 37
$ws.Range("E37").Value = 62

# Row 38
$ws.Range("E38").Value = 87
$ws.Range("F38").Value = 21
$ws.Range("H38").Value = 41

# Row 42
$ws.Range("E42").Value = 42

# Row 44
$ws.Range("E44").Value = 33

# Row 46
$ws.Range("E46").Value = 32
$ws.Range("F46").Value = 13
$ws.Range("H46").Value = 22

# Row 47
$ws.Range("E47").Value = 66

# Row 48
$ws.Range("E48").Value = 42

# Row 49
$ws.Range("E49").Value = 80
$ws.Range("F49").Value = 44
$ws.Range("H49").Value = 61

# Row 55
$ws.Range("E55").Value = 9

# Row 61
$ws.Range("E61").Value = 37
$ws.Range("F61").Value = 17
$ws.Range("H61").Value = 27

# Row 72
$ws.Range("E72").Value = 51

# Row 73
$ws.Range("E73").Value = 34
$ws.Range("F73").Value = 14
$ws.Range("H73").Value = 26

# Row 74
$ws.Range("E74").Value = 20
